$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update tag/status values: change "In-progress" to "Done" for rows 16 and 23 (col G)
$ws.Range("G16").Value = "Done"
$ws.Range("G23").Value = "Done"

# Update the view's top-left cell and active selection
$excel.ActiveWindow.ScrollRow = $ws.Range("B5").Row
$excel.ActiveWindow.ScrollColumn = $ws.Range("B5").Column
$ws.Range("G11").Select()
